# Update readme y archivo de estructura
$d = $word.ActiveDocument

# --- Locate the structure/naming-rules table (first table in the doc) ---
$t = $d.Tables(1)

# Resize the three columns (values are in points = twentieths-of-a-point/20)
$t.Columns(1).Width = 92.15   # 1843 dxa
$t.Columns(2).Width = 205.55  # 4111 dxa
$t.Columns(3).Width = 198.45  # 3969 dxa

# Overall preferred table width -> 9923 dxa (496.15 pt)
$t.PreferredWidth = 496.15

# Center the header-row text in all 3 header cells
$t.Cell(1,1).Range.ParagraphFormat.Alignment = 1
$t.Cell(1,2).Range.ParagraphFormat.Alignment = 1
$t.Cell(1,3).Range.ParagraphFormat.Alignment = 1

# --- Row 6 ("Trabajo Práctico Evaluable"): drop "trabajos_prácticos/" segment ---
$t.Cell(6,3).Range.Find.Execute("trabajos_prácticos/", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Row 7 ("Trabajo Práctico No Evaluable"): Practico/trabajos_practicos/ -> practico/ ---
$t.Cell(7,3).Range.Find.Execute("Practico/trabajos_practicos/", $true, $false, $false, $false, $false, $true, 1, $false, "practico/", 2) | Out-Null

# --- Bibliography list: merge the two split runs into one contiguous mention ---
$d.Content.Find.Execute("4K1_ISW_G3_Bibliografia_Agile_TESTING_-_A_Practical_Guide_For", $true, $false, $false, $false, $false, $true, 1, $false, "4K1_ISW_G3_Bibliografia_Agile_TESTING_-_A_Practical_Guide_For", 2) | Out-Null
